# Updated cryptos list (Price / Volume(1h)) - Tue Dec 26 03:19:24 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numeric-looking text (trailing zeros, thousands
# separators, etc. must be preserved verbatim as in the source feed), so force
# text formatting on the whole column before writing the new values, then put
# the number format back so no residual formatting is left on the cells.
$dPrices = $ws.Range("D2:D51")
$dPrices.NumberFormat = "@"

$ws.Range("D2").Value = "43.475.78"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "2.273.08"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "118.54"
$ws.Range("E5").Value = "  +5.76%  "
$ws.Range("D6").Value = "267.90"
$ws.Range("E6").Value = "  +1.32%  "
$ws.Range("D7").Value = "0.642"
$ws.Range("E7").Value = "  +3.52%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").Value = "0.621"
$ws.Range("E9").Value = "  +2.45%  "
$ws.Range("D10").Value = "47.28"
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("D11").Value = "0.0946"
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("D12").Value = "9.44"
$ws.Range("E12").Value = "  +7.67%  "
$ws.Range("E13").Value = "  -1.52%  "
$ws.Range("D14").Value = "15.78"
$ws.Range("E14").Value = "  +2.11%  "
$ws.Range("D15").Value = "0.917"
$ws.Range("E15").Value = "  +7.35%  "
$ws.Range("D16").Value = "2.613.62"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "2.266.59"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "43.679.88"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("D19").Value = "0.0000110"
$ws.Range("E19").Value = "  +1.66%  "
$ws.Range("D20").Value = "6.91"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").Value = "72.38"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("D22").Value = "2.40"
$ws.Range("E22").Value = "  -5.01%  "
$ws.Range("D23").Value = "234.68"
$ws.Range("E23").Value = "  +1.38%  "
$ws.Range("D24").Value = "2.97"
$ws.Range("E24").Value = "  +3.04%  "
$ws.Range("D25").Value = "9.69"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("D26").Value = "12.26"
$ws.Range("E26").Value = "  +8.63%  "
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("D28").Value = "41.86"
$ws.Range("E28").Value = "  +3.91%  "
$ws.Range("E29").Value = "  +2.04%  "
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("D31").Value = "174.53"
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("D32").Value = "21.55"
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("D33").Value = "0.0922"
$ws.Range("E33").Value = "  +2.15%  "
$ws.Range("D34").Value = "5.75"
$ws.Range("E34").Value = "  +0.98%  "
$ws.Range("D35").Value = "0.131"
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("D36").Value = "4.30"
$ws.Range("E36").Value = "  +12.82%  "
$ws.Range("D37").Value = "0.0382"
$ws.Range("E37").Value = "  +8.83%  "
$ws.Range("D38").Value = "4.63"
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("D39").Value = "0.108"
$ws.Range("D40").Value = "2.57"
$ws.Range("E40").Value = "  -0.80%  "
$ws.Range("D41").Value = "13.88"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").Value = "0.240"
$ws.Range("E42").Value = "  +2.19%  "
$ws.Range("D43").Value = "72.09"
$ws.Range("E43").Value = "  -5.28%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("D46").Value = "5.73"
$ws.Range("E46").Value = "  -5.65%  "
$ws.Range("D47").Value = "0.684"
$ws.Range("E47").Value = "  +22.33%  "
$ws.Range("D48").Value = "74.42"
$ws.Range("E48").Value = "  +37.94%  "
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("D50").Value = "103.50"
$ws.Range("E50").Value = "  +2.07%  "
$ws.Range("D51").Value = "8.59"
$ws.Range("E51").Value = "  -0.21%  "

$dPrices.Style = "Normal"
